# Plantilla e importación de horas extras
# Replace the human-readable Spanish column headers in row 1 with the
# internal Odoo field names (as documented in the "Guía" sheet), add a
# blank second row anchored at C2 (the new active selection), and strip
# the thin box border from the header row while centering it vertically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1: header labels -> Odoo technical field names ---------------
$ws.Range("A1").Value = "codigo_empleado"
$ws.Range("B1").Value = "periodo"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "horas_diurnas"
$ws.Range("E1").Value = "horas_nocturnas"
$ws.Range("F1").Value = "horas_diurnas_descanso"
$ws.Range("G1").Value = "horas_nocturnas_descanso"
$ws.Range("H1").Value = "horas_diurnas_asueto"
$ws.Range("I1").Value = "horas_nocturnas_asueto"

# --- Header row formatting: drop the box border, center vertically ----
$headerRow = $ws.Range("A1:I1")
$headerRow.Borders.LineStyle = -4142
$headerRow.VerticalAlignment = -4108

# --- New row 2: a single formatted (empty) cell at C2 ------------------
$ws.Range("C2").Font.Bold = $true

# --- Selection / active cell -------------------------------------------
$ws.Range("C2").Select()
